# Update cryptocurrency price/volume data per Sun Jan 29 07:50:23 UTC 2023 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value (both Price (D) and Volume(1h) (E) columns).
# NumberFormat is set to Text ("@") before assigning so Excel keeps these
# numeric-looking / percent-looking strings as literal text, matching the
# source data which is stored as text (inline strings) in the workbook.
$updates = @{
    "D2" = "308.59"
    "E2" = "-0.22%"
    "D3" = "38.97"
    "E3" = "-0.39%"
    "D4" = "5.122"
    "E4" = "0.33%"
    "D5" = "0.08120"
    "E5" = "-0.17%"
    "D6" = "1.945"
    "E6" = "-1.85%"
    "D7" = "8.140"
    "E7" = "2.69%"
    "D8" = "0.9267"
    "E8" = "-0.30%"
    "D9" = "0.1415"
    "E9" = "0.14%"
    "D10" = "0.1933"
    "E10" = "-1.19%"
    "D11" = "0.09021"
    "E11" = "-1.63%"
    "D12" = "0.03498"
    "E12" = "0.37%"
    "D13" = "0.09814"
    "E13" = "-0.21%"
    "D14" = "0.001403"
    "E14" = "-0.35%"
    "D15" = "0.006093"
    "E15" = "1.85%"
    "D16" = "3.906"
    "E16" = "8.68%"
    "D17" = "4.240"
    "E17" = "1.07%"
    "D18" = "3.357"
    "E18" = "-3.15%"
    "D19" = "0.3454"
    "E19" = "0.20%"
    "D20" = "0.1342"
    "E20" = "1.42%"
    "D21" = "4.737"
    "E21" = "-1.70%"
    "D22" = "0.2427"
    "E22" = "-1.70%"
    "D23" = "0.04378"
    "E23" = "-1.62%"
    "D24" = "0.001232"
    "E24" = "-0.47%"
    "D25" = "0.004795"
    "E25" = "-1.34%"
    "D26" = "0.0001302"
    "E26" = "-0.13%"
    "D27" = "0.0004006"
    "E27" = "-9.93%"
    "D39" = "0.02077"
    "E39" = "-1.43%"
    "D40" = "0.05105"
    "E40" = "-0.68%"
    "D41" = "0.007433"
    "E41" = "-0.63%"
    "D42" = "0.009801"
    "E42" = "-1.95%"
    "D43" = "0.1364"
    "E43" = "-0.09%"
    "D44" = "0.002133"
    "E44" = "-0.60%"
    "D45" = "0.008468"
    "E45" = "-16.58%"
    "D46" = "0.00006414"
    "E46" = "2.42%"
    "E47" = "-0.04%"
    "E48" = "-18.87%"
    "D49" = "0.002585"
    "D50" = "0.00002103"
    "E50" = "-0.04%"
    "D51" = "0.0002003"
    "E51" = "-0.04%"
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
}
